$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.766.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +5.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.706.58"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9985"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3681"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.53"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3300"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("E10").Value = "  +4.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07336"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9992"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.194"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.93"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.856"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.698.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.14%  "

$ws.Range("E17").Value = "  +3.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06630"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9982"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.055"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.70%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "25.741.54"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.468"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.483"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.58%  "

$ws.Range("E27").Value = "  +2.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.12"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.61%  "

$ws.Range("E29").Value = "  +10.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.888.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.06"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.89%  "

$ws.Range("E32").Value = "  +1.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.945"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08490"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.674"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.85"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.307"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.271"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06216"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.536"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2120"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02253"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.54"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +17.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6108"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9983"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("E46").Value = "  +2.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5832"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.89%  "

$ws.Range("E49").Value = "  +3.64%  "

$ws.Range("E50").Value = "  +4.99%  "

$ws.Range("E51").Value = "  +2.90%  "
